# Daily attendance processing - 2025-12-01 22:51:59
# Reorders the "Recorded By" (column G) values so that "System" is listed
# first among the recorders, for the two affected patterns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "System, backup@backdoor.com, system") {
        $cell.Value2 = "System, system, backup@backdoor.com"
    }
}
